# front panel: describe custom drawing
# Resize/reposition the "Rectangle 22" CCO label on Slide 1 and trim the
# library artifact name shown in it (drop the "-drawings" suffix).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item("Rectangle 22")   # id=23

# Reposition / resize (values chosen so the COM Single-precision round trip
# lands back on the exact target EMU: 335184 / 3622936 / 2767809).
$sh.Left = 26.39244194488189
$sh.Top = 285.27055418110234
$sh.Width = 217.93772153543307

# Update only the second run ("com.microej.clibrary.llimpl#microui-drawings"
# -> "com.microej.clibrary.llimpl#microui"), preserving the leading "CCO - "
# run and all run-level formatting.
$tr = $sh.TextFrame.TextRange
$old = "com.microej.clibrary.llimpl#microui-drawings"
$new = "com.microej.clibrary.llimpl#microui"
$start = $tr.Text.IndexOf($old) + 1
$chars = $tr.Characters($start, $old.Length)
$chars.Text = $new
